# Update cryptos price list (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.648.25"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "3.514.35"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("D7").Value = "3.514.22"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.432"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000217"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "4.106.11"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "3.511.17"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "67.601.14"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "449.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.635"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000129"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "3.651.42"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.04%  "
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").Value = "3.506.17"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0900"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "174.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "30.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.884"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("E48").Value = "  +3.56%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.87%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.255"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.31%  "
